$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.903
$ws.Range("C21").Value = -12.63
$ws.Range("C23").Value = -12.223
$ws.Range("C25").Value = -12.727
$ws.Range("C53").Value = -11.443
$ws.Range("C57").Value = -13.829
$ws.Range("C59").Value = -13.155
$ws.Range("C69").Value = -10.676
$ws.Range("C79").Value = -12.078
$ws.Range("C83").Value = -13.169
$ws.Range("C93").Value = -11.391
